$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "28.573.93"
$ws.Range('D3').Value = "1.566.64"
$ws.Range('E3').Value = "  -1.29%  "
$ws.Range('E4').Value = "  -0.06%  "
$ws.Range('D5').Value = "'211.83"
$ws.Range('E6').Value = "  -0.63%  "
$ws.Range('E7').Value = "  -0.05%  "
$ws.Range('D8').Value = "'46.39"
$ws.Range('E8').Value = "  +5.64%  "
$ws.Range('D9').Value = "'24.26"
$ws.Range('E9').Value = "  +1.63%  "
$ws.Range('E10').Value = "  -1.66%  "
$ws.Range('E11').Value = "  -1.36%  "
$ws.Range('D12').Value = "'0.0884"
$ws.Range('E12').Value = "  -0.39%  "
$ws.Range('D13').Value = "1.789.37"
$ws.Range('E13').Value = "  -1.35%  "
$ws.Range('D14').Value = "1.566.14"
$ws.Range('E14').Value = "  -1.39%  "
$ws.Range('E15').Value = "  -1.89%  "
$ws.Range('D16').Value = "28.570.76"
$ws.Range('E16').Value = "  +0.65%  "
$ws.Range('E17').Value = "  -2.95%  "
$ws.Range('D18').Value = "'62.15"
$ws.Range('E18').Value = "  -1.54%  "
$ws.Range('D19').Value = "'228.82"
$ws.Range('E19').Value = "  -1.38%  "
$ws.Range('D20').Value = "0.0₃0694"
$ws.Range('E20').Value = "  -2.11%  "
$ws.Range('E21').Value = "  -1.96%  "
$ws.Range('E22').Value = "  -0.14%  "
$ws.Range('E23').Value = "  -5.88%  "
$ws.Range('D24').Value = "'9.15"
$ws.Range('E24').Value = "  -2.68%  "
$ws.Range('E25').Value = "  +6.61%  "
$ws.Range('D26').Value = "'150.96"
$ws.Range('D27').Value = "'14.97"
$ws.Range('E27').Value = "  -1.80%  "
$ws.Range('E28').Value = "  -2.50%  "
$ws.Range('E29').Value = "  -3.36%  "
$ws.Range('E30').Value = "  -0.12%  "
$ws.Range('E31').Value = "  -1.39%  "
$ws.Range('E32').Value = "  -3.50%  "
$ws.Range('E33').Value = "  -0.98%  "
$ws.Range('D34').Value = "'3.12"
$ws.Range('E34').Value = "  -1.03%  "
$ws.Range('D35').Value = "1.398.95"
$ws.Range('E35').Value = "  -0.67%  "
$ws.Range('E36').Value = "  -0.76%  "
$ws.Range('E37').Value = "  -2.61%  "
$ws.Range('D38').Value = "'2.36"
$ws.Range('E38').Value = "  +0.83%  "
$ws.Range('D39').Value = "'2.58"
$ws.Range('E39').Value = "  +1.89%  "
$ws.Range('E40').Value = "  -0.94%  "
$ws.Range('D41').Value = "'0.538"
$ws.Range('E41').Value = "  -0.94%  "
$ws.Range('E43').Value = "  -3.41%  "
$ws.Range('E44').Value = "  +3.02%  "
$ws.Range('E45').Value = "  -4.30%  "
$ws.Range('E46').Value = "  -0.38%  "
$ws.Range('D47').Value = "'62.80"
$ws.Range('D48').Value = "1.702.14"
$ws.Range('E48').Value = "  -1.46%  "
$ws.Range('D49').Value = "'86.45"
$ws.Range('E49').Value = "  -1.40%  "
$ws.Range('E50').Value = "  -4.86%  "
$ws.Range('E51').Value = "  -1.04%  "
